# DESCW-1269 refine model & templates
#
# Re-works the report's column widths, the "data row" template font
# (8.5pt / bold-8pt -> a single 10pt font), drops all formatting from the
# blank spacer row (row 5), and updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------
# Project Name column becomes noticeably wider, the data columns get
# narrower to compensate.
$ws.Columns("B").ColumnWidth = 56
$ws.Columns("C:I").ColumnWidth = 10.833333333333332

# --- Template/data rows (3 & 4): font 8.5pt / bold 8pt -> 10pt -----
# Borders, number formats, fills and alignment are untouched - only the
# font size used throughout the two "{#r=...}" / "{#r1=...}" template
# rows changes.
$ws.Range("A3:I4").Font.Size = 10

# --- Row 5 (blank spacer row): strip all formatting -----------------
$ws.Range("A5:I5").Style = "Normal"
$ws.Rows(5).RowHeight = 19

# --- Restore the saved selection ------------------------------------
$ws.Range("C1:I1").Select()
